# Insert a new data row at row 315 (this pushes the existing rows 315-336
# down to 316-337, matching the rest of the table's weekly series).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(315).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A315").Value = 5
$ws.Range("B315").Value = "Macroferia Regional de Talca"
$ws.Range("C315").Value = "Maule"
$ws.Range("D315").Value = 44714
$ws.Range("E315").Value = 7
$ws.Range("F315").Value = 100112023
$ws.Range("G315").Value = "Brócoli"
$ws.Range("H315").Value = "Sin especificar"
$ws.Range("I315").Value = "Primera"
$ws.Range("J315").Value = 3000
$ws.Range("K315").Value = 1000
$ws.Range("L315").Value = 1000
$ws.Range("M315").Value = 1000
$ws.Range("N315").Value = "$/unidad"
$ws.Range("O315").Value = "Región del Maule"
$ws.Range("P315").Value = 1000
$ws.Range("Q315").Value = 1
$ws.Range("R315").Value = "Hortaliza"
